$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7708309292793274
$ws.Range("B1").Value = 1.203548312187195
$ws.Range("C1").Value = 4.399590492248535
$ws.Range("D1").Value = 4.038678646087646
$ws.Range("E1").Value = 1.416980504989624
